$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.528300762176514
$ws.Range("B1").Value = 2.386903047561646
$ws.Range("C1").Value = 5.318450927734375
$ws.Range("D1").Value = 1.471071243286133
$ws.Range("E1").Value = 0.7816293835639954
